$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 278
$ws.Range("I12").Value = 256.66666
$ws.Range("K12").Value = 256.66666
$ws.Range("M12").Value = -86.66665999999998
$ws.Range("H15").Value = 1044493.4
$ws.Range("I15").Value = 1044493.4
$ws.Range("K15").Value = 3133480.2
$ws.Range("M15").Value = -3133311.2
$ws.Range("H17").Value = 2094.7288
$ws.Range("J17").Value = 2094.7288
$ws.Range("L17").Value = 6284.1864
$ws.Range("N17").Value = -6620.1864
$ws.Range("H19").Value = 1044.9
$ws.Range("I19").Value = 957.1429000000001
$ws.Range("J19").Value = 1249.6666
$ws.Range("K19").Value = 957.1429000000001
$ws.Range("L19").Value = 1249.6666
$ws.Range("M19").Value = -782.1429000000001
$ws.Range("N19").Value = -1599.6666
$ws.Range("H31").Value = 10999
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = ""
$ws.Range("H32").Value = 3780.2
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = ""
$ws.Range("H33").Value = 2100.524
$ws.Range("I33").Value = 356.3846
$ws.Range("J33").Value = 4934.75
$ws.Range("K33").Value = 356.3846
$ws.Range("L33").Value = 4934.75
$ws.Range("M33").Value = -127.3846
$ws.Range("N33").Value = -5392.75
$ws.Range("H40").Value = 10913.833
$ws.Range("J40").Value = 13709.429
$ws.Range("L40").Value = 13709.429
$ws.Range("N40").Value = -14059.429
$ws.Range("H51").Value = 9512.643
$ws.Range("J51").Value = 10143.363
$ws.Range("L51").Value = 10143.363
$ws.Range("N51").Value = -11111.363
$ws.Range("H55").Value = 79.8
$ws.Range("I55").Value = 79.8
$ws.Range("K55").Value = 79.8
$ws.Range("M55").Value = 134.2
$ws.Range("H70").Value = 2673.5715
$ws.Range("J70").Value = 2434.0625
$ws.Range("L70").Value = 7302.1875
$ws.Range("N70").Value = -7842.1875
$ws.Range("H73").Value = 2673.5715
$ws.Range("J73").Value = 2434.0625
$ws.Range("L73").Value = 7302.1875
$ws.Range("N73").Value = -9174.1875
$ws.Range("H76").Value = 5949.4375
$ws.Range("I76").Value = 4771
$ws.Range("K76").Value = 4771
$ws.Range("M76").Value = -4456
$ws.Range("H79").Value = 5949.4375
$ws.Range("I79").Value = 4771
$ws.Range("K79").Value = 4771
$ws.Range("M79").Value = -3679
$ws.Range("H80").Value = 804.26666
$ws.Range("I80").Value = 999.8
$ws.Range("J80").Value = 706.5
$ws.Range("K80").Value = 2999.4
$ws.Range("L80").Value = 2119.5
$ws.Range("M80").Value = -2001.4
$ws.Range("N80").Value = -4115.5
$ws.Range("H83").Value = 804.26666
$ws.Range("I83").Value = 999.8
$ws.Range("J83").Value = 706.5
$ws.Range("K83").Value = 8998.199999999999
$ws.Range("L83").Value = 6358.5
$ws.Range("M83").Value = -4006.199999999999
$ws.Range("N83").Value = -16342.5
$ws.Range("H86").Value = 7442.222
$ws.Range("I86").Value = 7821.294
$ws.Range("J86").Value = 6797.8
$ws.Range("K86").Value = 7821.294
$ws.Range("L86").Value = 6797.8
$ws.Range("M86").Value = -6698.294
$ws.Range("N86").Value = -9043.799999999999
$ws.Range("H88").Value = 3390.4167
$ws.Range("I88").Value = 3789.8
$ws.Range("J88").Value = 3105.1428
$ws.Range("K88").Value = 3789.8
$ws.Range("L88").Value = 3105.1428
$ws.Range("M88").Value = -3383.8
$ws.Range("N88").Value = -3917.1428
$ws.Range("H89").Value = 7442.222
$ws.Range("I89").Value = 7821.294
$ws.Range("J89").Value = 6797.8
$ws.Range("K89").Value = 39106.47
$ws.Range("L89").Value = 33989
$ws.Range("M89").Value = -33490.47
$ws.Range("N89").Value = -45221
$ws.Range("H91").Value = 3390.4167
$ws.Range("I91").Value = 3789.8
$ws.Range("J91").Value = 3105.1428
$ws.Range("K91").Value = 3789.8
$ws.Range("L91").Value = 3105.1428
$ws.Range("M91").Value = -2385.8
$ws.Range("N91").Value = -5913.1428
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""
$ws.Range("H98").Value = 1905.8889
$ws.Range("I98").Value = 1382.4286
$ws.Range("J98").Value = 3738
$ws.Range("K98").Value = 1382.4286
$ws.Range("L98").Value = 3738
$ws.Range("M98").Value = 115.5714
$ws.Range("N98").Value = -6734
$ws.Range("H107").Value = 719.9
$ws.Range("I107").Value = 719.9
$ws.Range("K107").Value = 719.9
$ws.Range("M107").Value = 1200.1
$ws.Range("H111").Value = 92828.586
$ws.Range("I111").Value = 690
$ws.Range("K111").Value = 2070
$ws.Range("M111").Value = 997
$ws.Range("H118").Value = 4995
$ws.Range("I118").Value = 4995
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 14985
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -13328
$ws.Range("N118").Value = ""
$ws.Range("H122").Value = 1905.8889
$ws.Range("I122").Value = 1382.4286
$ws.Range("J122").Value = 3738
$ws.Range("K122").Value = 4147.2858
$ws.Range("L122").Value = 11214
$ws.Range("M122").Value = -1697.2858
$ws.Range("N122").Value = -16114
$ws.Range("H132").Value = 2772.5
$ws.Range("I132").Value = 2761.1904
$ws.Range("K132").Value = 8283.5712
$ws.Range("M132").Value = -5753.5712
$ws.Range("H135").Value = 589.0625
$ws.Range("I135").Value = 621.6667
$ws.Range("J135").Value = 100
$ws.Range("K135").Value = 5595.0003
$ws.Range("L135").Value = 900
$ws.Range("M135").Value = -3060.0003
$ws.Range("N135").Value = -5970
$ws.Range("H138").Value = 3052.766
$ws.Range("I138").Value = 2687.6938
$ws.Range("J138").Value = 3691.6428
$ws.Range("K138").Value = 8063.0814
$ws.Range("L138").Value = 11074.9284
$ws.Range("M138").Value = -2923.0814
$ws.Range("N138").Value = -21354.9284
$ws.Range("H141").Value = 4085.24
$ws.Range("I141").Value = 3296.05
$ws.Range("J141").Value = 7242
$ws.Range("K141").Value = 9888.150000000001
$ws.Range("L141").Value = 21726
$ws.Range("M141").Value = -4708.150000000001
$ws.Range("N141").Value = -32086

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 811.4666999999999
$ws.Range("I2").Value = 752.7826
$ws.Range("K2").Value = 752.7826
$ws.Range("M2").Value = -639.7826
$ws.Range("H32").Value = 6659.983
$ws.Range("I32").Value = 3008.9216
$ws.Range("K32").Value = 3008.9216
$ws.Range("M32").Value = -2721.9216
$ws.Range("H45").Value = 1645.9375
$ws.Range("I45").Value = 1084.909
$ws.Range("K45").Value = 1084.909
$ws.Range("M45").Value = -707.9090000000001
$ws.Range("H61").Value = 6640.2856
$ws.Range("I61").Value = 6013.6665
$ws.Range("J61").Value = 10400
$ws.Range("K61").Value = 6013.6665
$ws.Range("L61").Value = 10400
$ws.Range("M61").Value = -5801.6665
$ws.Range("N61").Value = -10824
$ws.Range("H74").Value = 4375.636
$ws.Range("I74").Value = 1902.5454
$ws.Range("J74").Value = 9321.817999999999
$ws.Range("K74").Value = 1902.5454
$ws.Range("L74").Value = 9321.817999999999
$ws.Range("M74").Value = -1028.5454
$ws.Range("N74").Value = -11069.818
$ws.Range("H77").Value = 4375.636
$ws.Range("I77").Value = 1902.5454
$ws.Range("J77").Value = 9321.817999999999
$ws.Range("K77").Value = 9512.726999999999
$ws.Range("L77").Value = 46609.09
$ws.Range("M77").Value = -5144.726999999999
$ws.Range("N77").Value = -55345.09
$ws.Range("H88").Value = 2968.1667
$ws.Range("I88").Value = 2096.1
$ws.Range("J88").Value = 4058.25
$ws.Range("K88").Value = 2096.1
$ws.Range("L88").Value = 4058.25
$ws.Range("M88").Value = -1690.1
$ws.Range("N88").Value = -4870.25
$ws.Range("H91").Value = 2968.1667
$ws.Range("I91").Value = 2096.1
$ws.Range("J91").Value = 4058.25
$ws.Range("K91").Value = 2096.1
$ws.Range("L91").Value = 4058.25
$ws.Range("M91").Value = -692.0999999999999
$ws.Range("N91").Value = -6866.25
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""
$ws.Range("H102").Value = 3703.5
$ws.Range("I102").Value = 2296.3333
$ws.Range("J102").Value = 5110.6665
$ws.Range("K102").Value = 2296.3333
$ws.Range("L102").Value = 5110.6665
$ws.Range("M102").Value = -674.3332999999998
$ws.Range("N102").Value = -8354.666499999999
$ws.Range("H110").Value = 1863.3829
$ws.Range("I110").Value = 1832.381
$ws.Range("K110").Value = 1832.381
$ws.Range("M110").Value = 212.6189999999999
$ws.Range("H116").Value = 811.4666999999999
$ws.Range("I116").Value = 752.7826
$ws.Range("K116").Value = 752.7826
$ws.Range("M116").Value = 1541.2174
$ws.Range("H132").Value = 5816.364
$ws.Range("I132").Value = 2215.3076
$ws.Range("K132").Value = 6645.9228
$ws.Range("M132").Value = -4115.9228
$ws.Range("H133").Value = 43729
$ws.Range("J133").Value = 43729
$ws.Range("L133").Value = 43729
$ws.Range("N133").Value = -48789
$ws.Range("H136").Value = 6640.2856
$ws.Range("I136").Value = 6013.6665
$ws.Range("J136").Value = 10400
$ws.Range("K136").Value = 18040.9995
$ws.Range("L136").Value = 31200
$ws.Range("M136").Value = -15490.9995
$ws.Range("N136").Value = -36300

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 811.4666999999999
$ws.Range("I3").Value = 752.7826
$ws.Range("K3").Value = 752.7826
$ws.Range("M3").Value = -638.7826
$ws.Range("H22").Value = 1532.3636
$ws.Range("I22").Value = 1238.375
$ws.Range("J22").Value = 2316.3333
$ws.Range("K22").Value = 1238.375
$ws.Range("L22").Value = 2316.3333
$ws.Range("M22").Value = -1065.375
$ws.Range("N22").Value = -2662.3333
$ws.Range("H105").Value = 2927.4546
$ws.Range("I105").Value = 2754.5293
$ws.Range("J105").Value = 3515.4
$ws.Range("K105").Value = 2754.5293
$ws.Range("L105").Value = 3515.4
$ws.Range("M105").Value = -1007.5293
$ws.Range("N105").Value = -7009.4
$ws.Range("H107").Value = 3289.138
$ws.Range("I107").Value = 3212.65
$ws.Range("J107").Value = 3459.111
$ws.Range("K107").Value = 3212.65
$ws.Range("L107").Value = 3459.111
$ws.Range("M107").Value = -1292.65
$ws.Range("N107").Value = -7299.111
$ws.Range("H134").Value = 3372.3462
$ws.Range("I134").Value = 1835.5454
$ws.Range("J134").Value = 11824.75
$ws.Range("K134").Value = 5506.6362
$ws.Range("L134").Value = 35474.25
$ws.Range("M134").Value = -2971.6362
$ws.Range("N134").Value = -40544.25

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 471.15686
$ws.Range("I7").Value = 500.78125
$ws.Range("J7").Value = 421.26315
$ws.Range("K7").Value = 500.78125
$ws.Range("L7").Value = 421.26315
$ws.Range("M7").Value = -387.78125
$ws.Range("N7").Value = -647.26315
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""
$ws.Range("H22").Value = 2778
$ws.Range("I22").Value = 1435
$ws.Range("K22").Value = 1435
$ws.Range("M22").Value = -1085
$ws.Range("H31").Value = 6726.9697
$ws.Range("I31").Value = 2867.8635
$ws.Range("K31").Value = 2867.8635
$ws.Range("M31").Value = -2572.8635
$ws.Range("H33").Value = 3000
$ws.Range("I33").Value = 3000
$ws.Range("K33").Value = 3000
$ws.Range("M33").Value = -2621
$ws.Range("H34").Value = 6726.9697
$ws.Range("I34").Value = 2867.8635
$ws.Range("K34").Value = 2867.8635
$ws.Range("M34").Value = -2665.8635
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""
$ws.Range("H58").Value = 4298.95
$ws.Range("I58").Value = 3166.5833
$ws.Range("J58").Value = 5997.5
$ws.Range("K58").Value = 3166.5833
$ws.Range("L58").Value = 5997.5
$ws.Range("M58").Value = -2963.5833
$ws.Range("N58").Value = -6403.5
$ws.Range("I93").Value = 10203.5
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 10203.5
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -8331.5
$ws.Range("N93").Value = -33744
$ws.Range("H94").Value = 4474.25
$ws.Range("I94").Value = 3826.25
$ws.Range("J94").Value = 4798.25
$ws.Range("K94").Value = 3826.25
$ws.Range("L94").Value = 4798.25
$ws.Range("M94").Value = -3375.25
$ws.Range("N94").Value = -5700.25
$ws.Range("H99").Value = 2449.5789
$ws.Range("I99").Value = 2460.6667
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 2460.6667
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = -962.6667000000002
$ws.Range("N99").Value = -5246
$ws.Range("H105").Value = 2789.6428
$ws.Range("I105").Value = 2388.8462
$ws.Range("J105").Value = 8000
$ws.Range("K105").Value = 2388.8462
$ws.Range("L105").Value = 8000
$ws.Range("M105").Value = -641.8462
$ws.Range("N105").Value = -11494
$ws.Range("H107").Value = 1380.6666
$ws.Range("I107").Value = 1170.2354
$ws.Range("K107").Value = 1170.2354
$ws.Range("M107").Value = 749.7646
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
$ws.Range("H126").Value = 2449.5789
$ws.Range("I126").Value = 2460.6667
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 7382.000100000001
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = -4912.000100000001
$ws.Range("N126").Value = -11690
$ws.Range("H132").Value = 5185.8823
$ws.Range("I132").Value = 3666.1538
$ws.Range("K132").Value = 10998.4614
$ws.Range("M132").Value = -8468.4614
$ws.Range("H134").Value = 5545.1665
$ws.Range("I134").Value = 5071.278
$ws.Range("J134").Value = 6966.8335
$ws.Range("K134").Value = 15213.834
$ws.Range("L134").Value = 20900.5005
$ws.Range("M134").Value = -12678.834
$ws.Range("N134").Value = -25970.5005
$ws.Range("H136").Value = 4298.95
$ws.Range("I136").Value = 3166.5833
$ws.Range("J136").Value = 5997.5
$ws.Range("K136").Value = 9499.749899999999
$ws.Range("L136").Value = 17992.5
$ws.Range("M136").Value = -6949.749899999999
$ws.Range("N136").Value = -23092.5

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1579.6
$ws.Range("I9").Value = 950
$ws.Range("K9").Value = 2850
$ws.Range("M9").Value = -2626
$ws.Range("H60").Value = 939.3333
$ws.Range("I60").Value = 159
$ws.Range("K60").Value = 477
$ws.Range("M60").Value = -226
$ws.Range("H80").Value = 4099.3335
$ws.Range("J80").Value = 3815.8333
$ws.Range("L80").Value = 11447.4999
$ws.Range("N80").Value = -13319.4999
$ws.Range("H83").Value = 4099.3335
$ws.Range("J83").Value = 3815.8333
$ws.Range("L83").Value = 34342.4997
$ws.Range("N83").Value = -43702.4997
$ws.Range("H114").Value = 467.14285
$ws.Range("I114").Value = 467.14285
$ws.Range("K114").Value = 1401.42855
$ws.Range("M114").Value = 1852.57145
$ws.Range("H129").Value = 15171451
$ws.Range("J129").Value = 33354534
$ws.Range("L129").Value = 100063602
$ws.Range("N129").Value = -100073602
$ws.Range("H133").Value = 3162.2307
$ws.Range("I133").Value = 3618.1667
$ws.Range("K133").Value = 10854.5001
$ws.Range("M133").Value = -5794.500100000001
$ws.Range("H134").Value = 1939.9231
$ws.Range("I134").Value = 1939.9231
$ws.Range("K134").Value = 5819.7693
$ws.Range("M134").Value = -749.7692999999999

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5326.8486
$ws.Range("I70").Value = 4406.769
$ws.Range("J70").Value = 5924.9
$ws.Range("K70").Value = 4406.769
$ws.Range("L70").Value = 5924.9
$ws.Range("M70").Value = -4136.769
$ws.Range("N70").Value = -6464.9
$ws.Range("H73").Value = 5326.8486
$ws.Range("I73").Value = 4406.769
$ws.Range("J73").Value = 5924.9
$ws.Range("K73").Value = 4406.769
$ws.Range("L73").Value = 5924.9
$ws.Range("M73").Value = -3470.769
$ws.Range("N73").Value = -7796.9
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = ""
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = ""
$ws.Range("H97").Value = 932
$ws.Range("I97").Value = 449
$ws.Range("J97").Value = 1334.5
$ws.Range("K97").Value = 449
$ws.Range("L97").Value = 1334.5
$ws.Range("M97").Value = 47
$ws.Range("N97").Value = -2326.5
$ws.Range("H102").Value = 3867.5
$ws.Range("I102").Value = 2553.4614
$ws.Range("J102").Value = 6307.857
$ws.Range("K102").Value = 2553.4614
$ws.Range("L102").Value = 6307.857
$ws.Range("M102").Value = -931.4614000000001
$ws.Range("N102").Value = -9551.857
$ws.Range("H107").Value = 715.46155
$ws.Range("I107").Value = 600.2
$ws.Range("J107").Value = 1099.6666
$ws.Range("K107").Value = 600.2
$ws.Range("L107").Value = 1099.6666
$ws.Range("M107").Value = 1319.8
$ws.Range("N107").Value = -4939.6666
$ws.Range("H113").Value = 7450
$ws.Range("I113").Value = 4900
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 4900
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -2730
$ws.Range("N113").Value = -14340
$ws.Range("H122").Value = 7203.727
$ws.Range("I122").Value = 1290.3334
$ws.Range("K122").Value = 3871.0002
$ws.Range("M122").Value = -1421.0002
$ws.Range("H132").Value = 4434.1665
$ws.Range("I132").Value = 2728.6538
$ws.Range("J132").Value = 6449.773
$ws.Range("K132").Value = 8185.9614
$ws.Range("L132").Value = 19349.319
$ws.Range("M132").Value = -5655.9614
$ws.Range("N132").Value = -24409.319

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 115718.7
$ws.Range("J7").Value = 19284.715
$ws.Range("L7").Value = 19284.715
$ws.Range("N7").Value = -19508.715
$ws.Range("H16").Value = 1599.1515
$ws.Range("I16").Value = 1455.68
$ws.Range("K16").Value = 1455.68
$ws.Range("M16").Value = -1285.68
$ws.Range("H33").Value = 14507.5
$ws.Range("I33").Value = 14507.5
$ws.Range("K33").Value = 14507.5
$ws.Range("M33").Value = -14217.5
$ws.Range("H46").Value = 2829.9
$ws.Range("J46").Value = 5478.778
$ws.Range("L46").Value = 5478.778
$ws.Range("N46").Value = -5854.778
$ws.Range("H61").Value = 3672.1875
$ws.Range("I61").Value = 3150.3333
$ws.Range("K61").Value = 3150.3333
$ws.Range("M61").Value = -2948.3333
$ws.Range("H82").Value = 1652.2
$ws.Range("I82").Value = 1272.5
$ws.Range("J82").Value = 2221.75
$ws.Range("K82").Value = 1272.5
$ws.Range("L82").Value = 2221.75
$ws.Range("M82").Value = -911.5
$ws.Range("N82").Value = -2943.75
$ws.Range("H85").Value = 1652.2
$ws.Range("I85").Value = 1272.5
$ws.Range("J85").Value = 2221.75
$ws.Range("K85").Value = 1272.5
$ws.Range("L85").Value = 2221.75
$ws.Range("M85").Value = -24.5
$ws.Range("N85").Value = -4717.75
$ws.Range("H105").Value = 48662.332
$ws.Range("J105").Value = 48662.332
$ws.Range("L105").Value = 48662.332
$ws.Range("N105").Value = -55650.332
$ws.Range("H113").Value = 3672.1875
$ws.Range("I113").Value = 3150.3333
$ws.Range("K113").Value = 3150.3333
$ws.Range("M113").Value = -980.3332999999998
$ws.Range("H126").Value = 115718.7
$ws.Range("J126").Value = 19284.715
$ws.Range("L126").Value = 57854.145
$ws.Range("N126").Value = -62794.145
$ws.Range("H132").Value = 5504.7856
$ws.Range("I132").Value = 3264.8572
$ws.Range("K132").Value = 9794.571599999999
$ws.Range("M132").Value = -7264.571599999999
$ws.Range("H136").Value = 5088.25
$ws.Range("I136").Value = 1906.1364
$ws.Range("J136").Value = 8977.5
$ws.Range("K136").Value = 5718.4092
$ws.Range("L136").Value = 26932.5
$ws.Range("M136").Value = -3168.4092
$ws.Range("N136").Value = -32032.5

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 32450
$ws.Range("I43").Value = 20000
$ws.Range("J43").Value = 44900
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 44900
$ws.Range("M43").Value = -19851
$ws.Range("N43").Value = -45198
$ws.Range("H64").Value = 85250.336
$ws.Range("J64").Value = 85250.336
$ws.Range("L64").Value = 85250.336
$ws.Range("N64").Value = -85746.336
$ws.Range("H67").Value = 85250.336
$ws.Range("J67").Value = 85250.336
$ws.Range("L67").Value = 85250.336
$ws.Range("N67").Value = -86966.336
$ws.Range("H81").Value = 4504.2
$ws.Range("I81").Value = 4604.875
$ws.Range("J81").Value = 4101.5
$ws.Range("K81").Value = 9209.75
$ws.Range("L81").Value = 8203
$ws.Range("M81").Value = -8148.75
$ws.Range("N81").Value = -10325
$ws.Range("H84").Value = 4504.2
$ws.Range("I84").Value = 4604.875
$ws.Range("J84").Value = 4101.5
$ws.Range("K84").Value = 46048.75
$ws.Range("L84").Value = 41015
$ws.Range("M84").Value = -40744.75
$ws.Range("N84").Value = -51623
$ws.Range("H100").Value = 1499
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1499
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2998
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = -4080
$ws.Range("H103").Value = 39999
$ws.Range("J103").Value = 39999
$ws.Range("L103").Value = 39999
$ws.Range("N103").Value = -42343
$ws.Range("H113").Value = 1026.2307
$ws.Range("I113").Value = 667.36365
$ws.Range("K113").Value = 2002.09095
$ws.Range("M113").Value = 167.90905
$ws.Range("H122").Value = 3192.125
$ws.Range("I122").Value = 3081.5
$ws.Range("J122").Value = 3413.375
$ws.Range("K122").Value = 9244.5
$ws.Range("L122").Value = 10240.125
$ws.Range("M122").Value = -6794.5
$ws.Range("N122").Value = -15140.125
$ws.Range("H126").Value = 1961.4736
$ws.Range("I126").Value = 1793.7858
$ws.Range("K126").Value = 5381.357400000001
$ws.Range("M126").Value = -2911.357400000001
$ws.Range("H132").Value = 2935.6924
$ws.Range("I132").Value = 2401.3333
$ws.Range("K132").Value = 7203.999899999999
$ws.Range("M132").Value = -4673.999899999999
$ws.Range("H136").Value = 3956.3928
$ws.Range("I136").Value = 3512.6667
$ws.Range("J136").Value = 4755.1
$ws.Range("K136").Value = 10538.0001
$ws.Range("L136").Value = 14265.3
$ws.Range("M136").Value = -7988.000100000001
$ws.Range("N136").Value = -19365.3
